$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[-, -, -, 'MEC-3B-M.A.Comp.CAD / CAM']"

$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "[-, -, -, 'MEC-3B-M.A.Comp.CAD / CAM']"

$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "[-, -, -, 'MEC-3B-M.A.Comp.CAD / CAM']"

$ws.Range("D6").Value = "-"

$ws.Range("D7").Value = "-"

$ws.Range("E8").Value = "['MEC-3B-M.A.Comp.CAD / CAM', -, -, -]"

$ws.Range("D10").Value = "['MEC-3A-Usin. CNC', -, 'MEC-3A-M.A.Comp.CAD / CAM', -]"

$ws.Range("D11").Value = "['MEC-3A-Usin. CNC', -, 'MEC-3A-M.A.Comp.CAD / CAM', -]"
$ws.Range("E11").Value = "-"

$ws.Range("D12").Value = "['MEC-3A-Usin. CNC', -, 'MEC-3A-M.A.Comp.CAD / CAM', -]"
$ws.Range("E12").Value = "-"

$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "-"

$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"

$ws.Range("D16").Value = "[-, 'MEC-3A-M.A.Comp.CAD / CAM', -, 'MEC-3A-Usin. CNC']"

$ws.Range("C18").Value = "[-, 'ELM-2NA-CAM', 'MEC-2NA-Usin. CNC', -]"
$ws.Range("D18").Value = "['MEC-2NA-Usin. CNC', 'MEC-2NB-CAD/CAM', -, -]"
$ws.Range("E18").Value = "['ELM-2NA-CAM', -, -, 'MEC-2NB-CAD/CAM']"
$ws.Range("F18").Value = "-"

$ws.Range("C19").Value = "[-, 'ELM-2NA-CAM', 'MEC-2NA-Usin. CNC', -]"
$ws.Range("D19").Value = "['MEC-2NA-Usin. CNC', 'MEC-2NB-CAD/CAM', -, -]"
$ws.Range("E19").Value = "['ELM-2NA-CAM', -, -, -]"
$ws.Range("F19").Value = "-"

$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "[-, 'MEC-2NB-CAD/CAM', -, -]"
$ws.Range("F20").Value = "-"

$ws.Range("C21").Value = "-"
$ws.Range("F21").Value = "-"
